# Rename the "Property1"/"Record" worksheets to the unified DataNode/DataTable
# naming scheme (commit: "unify the conception of DataNode, DataTable, Entity."),
# and make the DataTable sheet the active/selected tab, matching the saved
# workbook state in the target file.

$wb = $excel.ActiveWorkbook

$nodeSheet = $wb.Worksheets.Item("Property1")
$nodeSheet.Name = "DataNode"

$tableSheet = $wb.Worksheets.Item("Record")
$tableSheet.Name = "DataTable"

# Activate the DataTable sheet so it becomes the workbook's active tab,
# mirroring the tabSelected/activeTab move seen in the target file.
$tableSheet.Activate()
